$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three rows that held "Durga Puja", "Victory Day (Bangladesh)" and
# "Christmas Day" (old rows 4-6). Deleting shifts everything below up by 3,
# which is what moves the stray formatted cells on (old) rows 7 & 9 up to
# rows 4 & 6.
$ws.Rows("4:6").Delete()

# Update the remaining two calendar entries.
$ws.Range("A2").Value = "New Year Holiday"
$ws.Range("B2").Value = "1/1/2025"
$ws.Range("C2").Value = "1/1/2025"
$ws.Range("D2").Value = "holiday"

$ws.Range("A3").Value = "Team Meating"
$ws.Range("B3").Value = "2/1/2025"
$ws.Range("C3").Value = "2/1/2025"
$ws.Range("D3").Value = "event"
